$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    4  = 1.846
    5  = 2.035
    6  = 2.25
    7  = 1.858
    8  = 2.158
    9  = 1.803
    10 = 1.953
    11 = 2.305
    12 = 1.953
    13 = 1.908
    14 = 2.183
    15 = 2.491
    16 = 1.846
    17 = 1.846
    18 = 1.976
    19 = 1.82
    20 = 1.858
    21 = 2.035
    22 = 2.531
    23 = 1.668
}

foreach ($row in $values.Keys) {
    $ws.Range("AB$row").Value = $values[$row]
}
